$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)

# --- Fix 1: merge the two consecutive empty paragraphs (between the
# "...pills inside" line and the "Micro Switch" line) into a single
# empty paragraph, inside the "Rectangle 171" text box. ---
$rect = $s.Shapes.Item("Rectangle 171")
$tr = $rect.TextFrame.TextRange
$chars = $tr.Characters(191, 2)
$chars.Delete()

# --- Fix 2: move the picture up so it no longer covers the text
# behind it (was overlapping "Micro Switch" section). ---
$pic = $s.Shapes.Item("Picture 172")
$pic.Left = 595.1055118110236
$pic.Top = 60.2985
